# NN needs to be corrected
# Update computed columns B-F and H for rows 2-9 (column G / "Площадь (га)" is unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 1262.28466796875;  C = 0.947;              D = 0.9221000075340271; E = 1.420600056648254;  F = 0.8230999708175659; H = 0.7718 }
    3 = @{ B = 1135.424072265625; C = 0.9083;              D = 0.9111;             E = 1.069200038909912;  F = 0.8026999831199646; H = 0.6742 }
    4 = @{ B = 782.131591796875;  C = 0.9137;              D = 0.9114;             E = 1.040899991989136;  F = 0.8173999786376953; H = 0.6775 }
    5 = @{ B = 872.0767211914062; C = 0.9199000000000001;  D = 0.9207;             E = 1.004299998283386;  F = 0.8400999903678894; H = 0.7594 }
    6 = @{ B = 1147.35986328125;  C = 0.9077;              D = 0.9076;             E = 0.9861000180244446; F = 0.832099974155426;  H = 0.6438 }
    7 = @{ B = 902.7155151367188; C = 0.9091;              D = 0.9067000150680542; E = 1.003299951553345;  F = 0.8435999751091003; H = 0.6359 }
    8 = @{ B = 1020.027526855469; C = 0.914;               D = 0.9133;             E = 0.9890999794006348; F = 0.8457000255584717; H = 0.6936 }
    9 = @{ B = 7122.01953125;     C = 0.9177999999999999;  D = 0.9121;             E = 1.420600056648254;  F = 0.8026999831199646; H = 4.8562 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("H$row").Value = $vals.H
}
